# Update the tracked word count (B2) to its new value. All dependent
# formulas (B4, D9, E9, D11, E11) recalculate automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 18007

# Move the active selection to match the author's last cursor position.
[void]$ws.Range("F22").Select()
